$wb = $excel.ActiveWorkbook

# The edit happens on the "Repayment schedule" worksheet (sheet3.xml): a new,
# blank column is inserted before column N, pushing the old N/O/P ("Late",
# "Outstanding", "Disbursement") one slot to the right (O/P/Q), and the sheet
# becomes the active tab of the workbook.
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate()

# Insert a new, blank column at N, shifting the old N:P ("Late",
# "Outstanding", "Disbursement") one slot right to O:Q. The new column
# inherits column M's width (but not its bestFit flag), same as Excel does
# when a column is inserted immediately after a copy/insert-cells operation.
$mWidth = $ws.Columns("M:M").ColumnWidth
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = $mWidth

# Update the sheet's used-range dimension and selection to match the new
# layout.
$ws.Range("S6").Select() | Out-Null
